$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (pushes existing rows 52..160 down to 53..161)
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with a new weekly record. It mirrors the
# record that used to be at row 52 (same market/region/category/variety),
# but represents a different week: different date, volume, weighted price,
# origin and $/Kg price.
$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value = "Bíobío"
$ws.Cells.Item(52, 4).Value = 44498
$ws.Cells.Item(52, 5).Value = 8
$ws.Cells.Item(52, 6).Value = 100112017
$ws.Cells.Item(52, 7).Value = "Apio"
$ws.Cells.Item(52, 8).Value = "Americana (o)"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 450
$ws.Cells.Item(52, 11).Value = 8000
$ws.Cells.Item(52, 12).Value = 8500
$ws.Cells.Item(52, 13).Value = 8278
$ws.Cells.Item(52, 14).Value = "$/docena de matas"
$ws.Cells.Item(52, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 16).Value = 1380
$ws.Cells.Item(52, 17).Value = 6
$ws.Cells.Item(52, 18).Value = "Hortaliza"
